$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 191 — this shifts the existing rows 191..224 down to 192..225
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new data record
$ws.Cells.Item(191, 1).Value2  = 5
$ws.Cells.Item(191, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(191, 3).Value2  = "Maule"
$ws.Cells.Item(191, 4).Value2  = 44504
$ws.Cells.Item(191, 4).NumberFormat = $ws.Cells.Item(192, 4).NumberFormat
$ws.Cells.Item(191, 5).Value2  = 7
$ws.Cells.Item(191, 6).Value2  = 100112032
$ws.Cells.Item(191, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(191, 8).Value2  = "Sin especificar"
$ws.Cells.Item(191, 9).Value2  = "Primera"
$ws.Cells.Item(191, 10).Value2 = 400
$ws.Cells.Item(191, 11).Value2 = 7000
$ws.Cells.Item(191, 12).Value2 = 7000
$ws.Cells.Item(191, 13).Value2 = 7000
$ws.Cells.Item(191, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(191, 15).Value2 = "Región del Maule"
$ws.Cells.Item(191, 16).Value2 = 117
$ws.Cells.Item(191, 17).Value2 = 60
$ws.Cells.Item(191, 18).Value2 = "Hortaliza"
